$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing A3 timestamp (tiny precision correction)
$ws.Cells.Item(3, 1).Value = 45878.08351149305

# Add new row 4 with data
$ws.Cells.Item(4, 1).Value = 45878.12517276972
$ws.Cells.Item(4, 1).NumberFormat = $ws.Cells.Item(3, 1).NumberFormat

$ws.Cells.Item(4, 2).Value = 2025
$ws.Cells.Item(4, 3).Value = 37
$ws.Cells.Item(4, 4).Value = 13.28
$ws.Cells.Item(4, 5).Value = 91.91
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 2.3
$ws.Cells.Item(4, 8).Value = "NNW"
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = "03:00:14"
